# ---------------------------------------------------------------------------
# Helper: replace a literal substring inside a shape's text, preserving the
# rest of the TextFrame (other paragraphs / runs) untouched. Returns $true
# if the substring was found & replaced.
# ---------------------------------------------------------------------------
function Replace-InShapeText {
    param(
        $Shape,
        [string]$OldText,
        [string]$NewText
    )
    $tr = $Shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -lt 0) {
        return $false
    }
    $sub = $tr.Characters($idx + 1, $OldText.Length)
    $sub.Text = $NewText
    return $true
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer "date" placeholder on every slide layout: 30-05-2022 -> 08-06-2022
#    (the presentation was re-saved on 2022-06-08; the automatically-updating
#    datetimeFigureOut field caches the new date text).
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            Replace-InShapeText $sh "30-05-2022" "08-06-2022" | Out-Null
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 1 - subtitle "May 2022": split "May" off into its own run so it
#    reverts to the placeholder's default (18pt) size while " 2022" keeps
#    the explicit 24pt size it already had.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$subtitle1 = $slide1.Shapes.Item(2)
$tr1 = $subtitle1.TextFrame.TextRange
$full1 = $tr1.Text
$mayIdx = $full1.IndexOf("May 2022")
if ($mayIdx -ge 0) {
    $maySub = $tr1.Characters($mayIdx + 1, 3)
    $maySub.Font.Size = 18
}

# ---------------------------------------------------------------------------
# 3) Slide 4 - "Solution Approach" bullet list text tweaks.
# ---------------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$subtitle4 = $slide4.Shapes.Item(2)

Replace-InShapeText $subtitle4 `
    "We have used Azure Logic Apps and Azure ML to implement the solution." `
    "Azure Logic Apps and Azure ML to implement the solution." | Out-Null

Replace-InShapeText $subtitle4 `
    "Used Low/Code, No Code Approach." `
    "Low Code, No Code Approach." | Out-Null

Replace-InShapeText $subtitle4 `
    "Most implementation done with Azure Logic Apps(serverless), that simplifies deployment" `
    "Most of the implementation done with Azure Logic Apps (serverless), that simplifies deployment." | Out-Null

# ---------------------------------------------------------------------------
# 4) Slide 10 - remove stray double-space before the closing parenthesis.
# ---------------------------------------------------------------------------
$slide10 = $p.Slides.Item(10)
$subtitle10 = $slide10.Shapes.Item(2)

Replace-InShapeText $subtitle10 `
    "Give a diversity score based on the number of diversity parameters and their priority category met (e.g. of priority: Woman +African American>  White + Differently Abled )" `
    "Give a diversity score based on the number of diversity parameters and their priority category met (e.g. of priority: Woman +African American>  White + Differently Abled)" | Out-Null
